# Update gh-pages generated data (matches commit "Update gh-pages to output generated at 456a3b4")
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 151
$ws1.Range("F4").Value  = 62
$ws1.Range("F5").Value  = 502
$ws1.Range("F6").Value  = 1487
$ws1.Range("F7").Value  = 918
$ws1.Range("F8").Value  = 109
$ws1.Range("F9").Value  = 206
$ws1.Range("F10").Value = 147
$ws1.Range("F11").Value = 201
$ws1.Range("F12").Value = 119
$ws1.Range("G12").Value = 45
$ws1.Range("F13").Value = 173
$ws1.Range("F14").Value = 158

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# --- Sheet "全部类型" (all types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 151
$ws4.Range("F4").Value  = 62
$ws4.Range("F5").Value  = 502
$ws4.Range("F6").Value  = 1487
$ws4.Range("F7").Value  = 2
$ws4.Range("F8").Value  = 918
$ws4.Range("F9").Value  = 109
$ws4.Range("F10").Value = 206
$ws4.Range("F11").Value = 147
$ws4.Range("F12").Value = 201
$ws4.Range("F13").Value = 119
$ws4.Range("G13").Value = 45
$ws4.Range("F14").Value = 173
$ws4.Range("F15").Value = 158
